# "Changement camera et decor"
#
# - Fix a typo in the "lifebar design" task label (A4)
# - Add a new backlog task "Coup corps à corps" (A31)
# - Extend the table with additional empty (but bordered) rows so the
#   sheet keeps room for more tasks (rows 32-41)
# - Leave the selection on A17, matching the saved view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo: "Designe" -> "Design"
$ws.Range("A4").Value = "Design Barre de vie : style mortal kombat X"

# New task row
$ws.Range("A31").Value = "Coup corps à corps"

# Apply the same thin-border look used by the rest of the table to the
# newly added rows (31 through 41, columns A-D)
$ws.Range("A31:D41").Borders.LineStyle = 1

# Restore the saved selection/active cell
$ws.Range("A17").Select() | Out-Null
